$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "The moving platform consists of the base waterjet cut piece",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The moving platform consists of the base waterjet cut piece", 2
)
